# The deck's Slide Master theme (Integral) is swapped for the Office
# default theme (Office Theme) -- i.e. the presentation's active colour
# scheme changes from the "Integral" palette to the standard "Office"
# palette. Font scheme and format scheme are already identical between
# the old and new theme, so only the twelve theme colours actually move.
#
# PowerPoint's ThemeColor.RGB takes a COM COLORREF (0x00BBGGRR), i.e.
# R + G*256 + B*65536.
function ColorRefFromHex([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Theme colour scheme slots, in PowerPoint's fixed Item() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeTheme = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$cs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $cs.Item($i).RGB = ColorRefFromHex $officeTheme[$i - 1]
}
